# Refresh the "cryptos" price/volume snapshot (GitHub Actions data pull).
# Columns: D = Price (text, dotted-thousands formatting from the source feed),
#          E = Volume(1h) change (text, padded "  +x.xx%  " strings).
# A handful of Price cells are plain-number-looking strings (e.g. "211.55",
# "7.00"); Excel's COM layer auto-coerces those to numeric on a bare
# Value= assignment (and "7.00" would lose its trailing zero as 7). To keep
# them stored as text -- matching the original inline-string cells, with no
# residual cell formatting -- force NumberFormat to "@" before the write and
# restore the default "Normal" style immediately afterward.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.669.74'
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("D3").Value = '1.597.92'
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = '  +0.92%  '
$ws.Range("E7").Value = '  +0.19%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  +0.66%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.55'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0838'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.09%  '
$ws.Range("D12").Value = '1.822.21'
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").Value = '1.609.74'
$ws.Range("E13").Value = '  +0.73%  '
$ws.Range("E14").Value = '  -0.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.523'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.31%  '
$ws.Range("E16").Value = '  +0.13%  '
$ws.Range("D17").Value = '26.652.75'
$ws.Range("D18").Value = '0.0₃0736'
$ws.Range("E18").Value = '  +1.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '209.58'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("E20").Value = '  +0.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.45%  '
$ws.Range("E22").Value = '  +0.48%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.34'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.63%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.99'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.45'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.50%  '
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("E27").Value = '  -0.64%  '
$ws.Range("E28").Value = '  -0.51%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.28'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("E30").Value = '  +2.43%  '
$ws.Range("E31").Value = '  +0.24%  '
$ws.Range("E32").Value = '  +0.95%  '
$ws.Range("E33").Value = '  +1.56%  '
$ws.Range("D34").Value = '1.285.62'
$ws.Range("E34").Value = '  -1.00%  '
$ws.Range("E35").Value = '  -6.77%  '
$ws.Range("E36").Value = '  +0.57%  '
$ws.Range("E37").Value = '  +0.91%  '
$ws.Range("E38").Value = '  -0.67%  '
$ws.Range("E39").Value = '  -1.01%  '
$ws.Range("E40").Value = '  +20.03%  '
$ws.Range("E41").Value = '  +2.18%  '
$ws.Range("E42").Value = '  -0.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.785'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.61%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.64'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.33%  '
$ws.Range("D45").Value = '1.735.36'
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '90.64'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.47%  '
$ws.Range("E47").Value = '  -3.44%  '
$ws.Range("D48").Value = '0.0₆0105'
$ws.Range("E48").Value = '  -0.12%  '
$ws.Range("E49").Value = '  +0.87%  '
$ws.Range("E50").Value = '  +0.80%  '
$ws.Range("E51").Value = '  +0.02%  '
